$d = $word.ActiveDocument

# --- First paragraph: ID placeholder text + formatting updates ---
$p1 = $d.Paragraphs(1)

$oldText = "**ID__AFFARS_pgi_5301_topic_40__ID**"
$newText = "**ID__AFFARS_SMC_PGI_5301__ID**"

# Replace just the first run's text (leave the following space/run alone for now)
$start = $p1.Range.Start
$idRange = $d.Range($start, $start + $oldText.Length)
$idRange.Text = $newText

# Remove the now-orphaned trailing space that used to sit in the 2nd run
$spaceStart = $start + $newText.Length
$spaceRange = $d.Range($spaceStart, $spaceStart + 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

# Paragraph indent: 120 twips (6pt) -> 225 twips (11.25pt)
$p1.Format.LeftIndent = 11.25

# Paragraph border on all four sides, space="5" (no line drawn, just the
# reserved spacing), matching <w:pBdr><w:top w:space="5"/>...</w:pBdr>
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
